# Sincronizando os dados entre documentos
# - Remove os itens de "ingresso disponível" e renumera os RF's (RF-14..RF-16
#   passam a ser usados para o CRUD de "Categoria de Evento", e o restante da
#   numeração é deslocado).
# - "Cadastrar/Exibir/Editar/Excluir ingresso disponível" passam a se chamar
#   "Cadastrar/Exibir/Editar/Excluir tipo de ingresso".
# - Seleção da planilha passa a ser a faixa inteira de dados (A1:C31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renumeração da coluna "Identificador" (coluna A) para as linhas 15 a 31.
$ws.Cells.Item(15, 1).Value = "RF-14"
$ws.Cells.Item(16, 1).Value = "RF-15"
$ws.Cells.Item(17, 1).Value = "RF-16"
$ws.Cells.Item(18, 1).Value = "RF-17"
$ws.Cells.Item(19, 1).Value = "RF-18"
$ws.Cells.Item(20, 1).Value = "RF-19"
$ws.Cells.Item(21, 1).Value = "RF-20"
$ws.Cells.Item(22, 1).Value = "RF-21"
$ws.Cells.Item(23, 1).Value = "RF-22"
$ws.Cells.Item(24, 1).Value = "RF-23"
$ws.Cells.Item(25, 1).Value = "RF-24"
$ws.Cells.Item(26, 1).Value = "RF-25"
$ws.Cells.Item(27, 1).Value = "RF-26"
$ws.Cells.Item(28, 1).Value = "RF-27"
$ws.Cells.Item(29, 1).Value = "RF-28"
$ws.Cells.Item(30, 1).Value = "RF-29"
$ws.Cells.Item(31, 1).Value = "RF-30"

# Renomeação da "Descrição Resumida" (coluna B) para as linhas 24 a 27:
# "ingresso disponível" -> "tipo de ingresso".
$ws.Cells.Item(24, 2).Value = "Cadastrar tipo de ingresso"
$ws.Cells.Item(25, 2).Value = "Exibir tipo de ingresso"
$ws.Cells.Item(26, 2).Value = "Editar tipo de ingresso"
$ws.Cells.Item(27, 2).Value = "Excluir tipo de ingresso"

# Atualiza a seleção ativa da planilha para cobrir toda a faixa de dados.
$ws.Range("A1:C31").Select()
